$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01675466666666667
$ws.Range("H2").Value = 0.050264
$ws.Range("I2").Value = 0.0001854906931657378
$ws.Range("J2").Value = 0.0001854906931657378
$ws.Range("M2").Value = 1.343359
$ws.Range("N2").Value = 4.030077
$ws.Range("O2").Value = 0.736296379391111
$ws.Range("P2").Value = 0.7362963793911109
$ws.Range("Q2").Value = 0.02250753225866667
$ws.Range("R2").Value = 0.202567790328
$ws.Range("S2").Value = 0.0001365761257886802
$ws.Range("T2").Value = 0.0001365761257886802
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01675466666666667
$ws.Range("H3").Value = 0.050264
$ws.Range("I3").Value = 0.0001854906931657378
$ws.Range("J3").Value = 0.0001854906931657378
$ws.Range("O3").Value = 0.1764523396969075
$ws.Range("P3").Value = 0.1764523396969075
$ws.Range("Q3").Value = 0.005393896858666667
$ws.Range("R3").Value = 0.048545071728
$ws.Range("S3").Value = 0.000032730266801095586831973905
$ws.Range("T3").Value = 0.000032730266801095586831973905
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01675466666666667
$ws.Range("H4").Value = 0.050264
$ws.Range("I4").Value = 0.0001854906931657378
$ws.Range("J4").Value = 0.0001854906931657378
$ws.Range("O4").Value = 0.08725128091198156
$ws.Range("P4").Value = 0.08725128091198156
$ws.Range("Q4").Value = 0.002667147462222223
$ws.Range("R4").Value = 0.02400432716
$ws.Range("S4").Value = 0.000016184300575961970515082369
$ws.Range("T4").Value = 0.000016184300575961960350687002
$ws.Range("I5").Value = 0.9933938536206305
$ws.Range("J5").Value = 0.9933938536206304
$ws.Range("M5").Value = 1.343359
$ws.Range("N5").Value = 4.030077
$ws.Range("O5").Value = 0.736296379391111
$ws.Range("P5").Value = 0.7362963793911109
$ws.Range("Q5").Value = 120.538900493243
$ws.Range("R5").Value = 1084.850104439187
$ws.Range("S5").Value = 0.7314322977302535
$ws.Range("T5").Value = 0.7314322977302534
$ws.Range("I6").Value = 0.9933938536206305
$ws.Range("J6").Value = 0.9933938536206304
$ws.Range("O6").Value = 0.1764523396969075
$ws.Range("P6").Value = 0.1764523396969075
$ws.Range("S6").Value = 0.1752866697118874
$ws.Range("T6").Value = 0.1752866697118874
$ws.Range("I7").Value = 0.9933938536206305
$ws.Range("J7").Value = 0.9933938536206304
$ws.Range("O7").Value = 0.08725128091198156
$ws.Range("P7").Value = 0.08725128091198156
$ws.Range("S7").Value = 0.08667488617848952
$ws.Range("T7").Value = 0.08667488617848951
$ws.Range("G8").Value = 0.5799533333333334
$ws.Range("I8").Value = 0.006420655686203657
$ws.Range("J8").Value = 0.006420655686203655
$ws.Range("M8").Value = 1.343359
$ws.Range("N8").Value = 4.030077
$ws.Range("O8").Value = 0.736296379391111
$ws.Range("P8").Value = 0.7362963793911109
$ws.Range("Q8").Value = 0.7790855299133336
$ws.Range("R8").Value = 7.011769769220002
$ws.Range("S8").Value = 0.004727505535068702
$ws.Range("T8").Value = 0.0047275055350687
$ws.Range("G9").Value = 0.5799533333333334
$ws.Range("I9").Value = 0.006420655686203657
$ws.Range("J9").Value = 0.006420655686203655
$ws.Range("O9").Value = 0.1764523396969075
$ws.Range("P9").Value = 0.1764523396969075
$ws.Range("Q9").Value = 0.1867066964133334
$ws.Range("S9").Value = 0.001132939718218888
$ws.Range("T9").Value = 0.001132939718218888
$ws.Range("G10").Value = 0.5799533333333334
$ws.Range("I10").Value = 0.006420655686203657
$ws.Range("J10").Value = 0.006420655686203655
$ws.Range("O10").Value = 0.08725128091198156
$ws.Range("P10").Value = 0.08725128091198156
$ws.Range("Q10").Value = 0.09232180454444447
$ws.Range("R10").Value = 0.8308962409000001
$ws.Range("S10").Value = 0.000560210432916067
$ws.Range("T10").Value = 0.0005602104329160668
